$wb = $excel.ActiveWorkbook

# The edited sheet is "Nädal 4" (the 4th week sheet, last tab).
$ws = $wb.Worksheets.Item("Nädal 4")

# Column G (Activity) needs to widen to fit the new, longer entry text.
$ws.Columns.Item(7).ColumnWidth = 20.6640625

# --- Row 13 (entry #7): new time-log entry -------------------------------
$ws.Cells.Item(13, 2).Value = 43884                        # B13 date
$ws.Cells.Item(13, 3).Value = 0.47569444444444442          # C13 start time
$ws.Cells.Item(13, 4).Value = 0.60069444444444442          # D13 stop time
$ws.Cells.Item(13, 5).Value = 40                            # E13 interruption (min)
$ws.Cells.Item(13, 6).Value = 140                           # F13 delta time (min)
$ws.Cells.Item(13, 7).Value = "Pluralsight - videokursus"   # G13 activity
$ws.Cells.Item(13, 8).Value = "Razor Pages in ASP.NET Core: Getting Started (part 1-3)" # H13 comments
$ws.Cells.Item(13, 10).Value = "x"                          # J13 marked "U"

# --- Row 14 (entry #8): start of next entry -------------------------------
$ws.Cells.Item(14, 2).Value = 43884                         # B14 date
$ws.Cells.Item(14, 3).Value = 0.6069444444444444            # C14 start time
$ws.Cells.Item(14, 7).Value = "Kodutöö 4"                   # G14 activity

# --- Selection moves to G14, matching the saved view ----------------------
$ws.Range("G14").Select()
